$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update product_type (column F) for rows 15-17 from "gamebook" to "slipcase set"
$ws.Range("F15:F17").Value = "slipcase set"

# Update the active selection to match the edited range
$ws.Range("F15:F17").Select()
